$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.003.73"
$ws.Range("E2").Value = "  -2.21%  "
$ws.Range("D3").Value = "1.666.86"
$ws.Range("E3").Value = "  -1.95%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.67"
$ws.Range("E5").Value = "  -1.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5095"
$ws.Range("E6").Value = "  -0.50%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2641"
$ws.Range("E8").Value = "  -0.57%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06390"
$ws.Range("E9").Value = "  +1.47%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.93"
$ws.Range("E10").Value = "  -1.42%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07412"
$ws.Range("E11").Value = "  +0.77%  "
$ws.Range("D12").Value = "1.663.78"
$ws.Range("E12").Value = "  -2.23%  "
$ws.Range("E13").Value = "  -0.66%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5824"
$ws.Range("E14").Value = "  -0.41%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008494"
$ws.Range("E15").Value = "  +0.25%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.17"
$ws.Range("E16").Value = "  -2.49%  "
$ws.Range("D17").Value = "26.064.24"
$ws.Range("E17").Value = "  -2.18%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.934"
$ws.Range("E18").Value = "  -2.08%  "
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("E20").Value = "  -2.53%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "189.85"
$ws.Range("E21").Value = "  +1.35%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.202"
$ws.Range("E22").Value = "  -1.47%  "
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "145.01"
$ws.Range("E24").Value = "  +0.32%  "
$ws.Range("E25").Value = "  +0.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1194"
$ws.Range("E26").Value = "  +3.32%  "
$ws.Range("E27").Value = "  -0.57%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06718"
$ws.Range("E28").Value = "  +17.28%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.312"
$ws.Range("E29").Value = "  -1.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.313"
$ws.Range("E30").Value = "  -2.26%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.518"
$ws.Range("E31").Value = "  -0.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.508"
$ws.Range("E32").Value = "  -0.10%  "
$ws.Range("E33").Value = "  -1.16%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.017"
$ws.Range("E34").Value = "  -1.34%  "
$ws.Range("E35").Value = "  +0.51%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.368"
$ws.Range("E36").Value = "  -0.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.713"
$ws.Range("E37").Value = "  +1.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.212"
$ws.Range("E38").Value = "  +6.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01600"
$ws.Range("E39").Value = "  -0.77%  "
$ws.Range("D40").Value = "1.075.71"
$ws.Range("E40").Value = "  -2.68%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8592"
$ws.Range("E41").Value = "  -0.36%  "
$ws.Range("E42").Value = "  +0.47%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.54"
$ws.Range("E43").Value = "  +1.34%  "
$ws.Range("D44").Value = "1.814.45"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000115"
$ws.Range("E45").Value = "  +3.45%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.26"
$ws.Range("E46").Value = "  -1.14%  "
$ws.Range("E47").Value = "  +0.39%  "
$ws.Range("E48").Value = "  -1.82%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05209"
$ws.Range("E49").Value = "  -0.80%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4289"
$ws.Range("E50").Value = "  -0.89%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.953"
$ws.Range("E51").Value = "  +2.21%  "
